# "Adding class names to all of the classes" -- fills in the GameView
# class-writeup placeholders with the real class name.
$d = $word.ActiveDocument

# 1) Title heading: "Some Class Name" -> "GameView"
$d.Content.Find.Execute("Some Class Name", $true, $false, $false, $false, $false,
                         $true, 1, $false, "GameView", 2) | Out-Null

# 2) CRC table header cell: "Class Name" -> "GameView"
#    ("Some Class Name" is already gone, so this is now unambiguous.)
$d.Content.Find.Execute("Class Name", $true, $false, $false, $false, $false,
                         $true, 1, $false, "GameView", 2) | Out-Null

# 3) The description paragraph had been mid-edit: the word "Some" was split
#    into "So" + "me" around the (cursor-position) _GoBack bookmark. Collapse
#    that back into a single, normal "Some paragraph about this class" run,
#    using the bookmark itself to find the exact split point so we don't
#    touch any of the many other identical "Some paragraph..." sentences.
$goBack = $d.Bookmarks("_GoBack")
$splitRange = $d.Range($goBack.Start - 2, $goBack.End + 29)
$splitRange.Text = "Some paragraph about this class"
$d.Bookmarks("_GoBack").Delete()

# 4) Re-create the _GoBack bookmark where Word would actually leave it: right
#    after the class name that was just typed into the CRC table cell.
$tbl = $d.Tables(1)
$cellRange = $tbl.Cell(1, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null   # wdCharacter; trim the cell-end mark
$d.Bookmarks.Add("_GoBack", $cellRange) | Out-Null
